# Update the player_elo_table sheet with the new game results
# ("beach week" and "the last roll" games): new players added,
# several Elo scores updated, and the row ordering re-sorted by Elo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data set (player name, current_elo), row 2 .. row 34
$data = @(
    @("Aaron Carter", 1578.068408583175),
    @("Will Simpson", 1389.845978659066),
    @("Kevin Lee", 1385.601170521426),
    @("Roman Ramirez", 1354.964499748664),
    @("Rohan Chowla", 1303.641964193648),
    @("Kevin Cooper", 1271.813139769939),
    @("Bronte Sundstrom", 1244.912353653355),
    @("Jack Massingill", 1228.341892966386),
    @("Nathan Snow", 1219.740685506919),
    @("Kristian Banlaoi", 1216.037581761179),
    @("Yvonne Nguyen", 1204.868130393788),
    @("Gabe Silverstein", 1201.692490161548),
    @("Piper Parker", 1194.748819295017),
    @("Parker Simpson", 1193.5144450266),
    @("Carla Betancourt", 1192.258245275334),
    @("Cason Duszak", 1189.056587063621),
    @("Reagan Fryatt", 1182.605339128595),
    @("Spencer Harris", 1179.730047956958),
    @("Abby LeahFriend", 1178.180583312953),
    @("Anna Brown", 1176.279801661901),
    @("Leah Baetcke", 1167.526278792758),
    @("Coby Lovelace", 1160.048928653962),
    @("Rose Roché", 1145.11049623546),
    @("Helen Dunn", 1143.988590686295),
    @("Noah Dale", 1143.988590686295),
    @("Ann Hall", 1142.733815855265),
    @("Jason Jackson", 1140.220115271984),
    @("Brian Tafazoli", 1140),
    @("Matthew Rusten", 1139.327236294392),
    @("Evan Sooklal", 1113.066776939508),
    @("Cassie Deering", 1058.266119428515),
    @("Sam Tellis", 1040.370685155026),
    @("Paul Bartenfeld", 979.4502013604686)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row = $row + 1
}
